$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update forecast coefficient values
$ws.Range("B2").Value = 0.5588025657981637
$ws.Range("B3").Value = -0.2901569498637481
$ws.Range("B4").Value = 3.351163668276367

# Remove row 5 entirely (A5="4", B5=0.3218894107672925)
$ws.Range("A5:B5").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
